$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.465929666666667
$ws.Range("H2").Value = 4.397789
$ws.Range("M2").Value = 1.660421
$ws.Range("N2").Value = 4.981262999999999
$ws.Range("O2").Value = 0.03714789785507311
$ws.Range("P2").Value = 0.03714789785507311
$ws.Range("Q2").Value = 2.434060403056333
$ws.Range("R2").Value = 21.90654362750699
$ws.Range("S2").Value = 0.03714789785507311
$ws.Range("T2").Value = 0.03714789785507311

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.465929666666667
$ws.Range("H3").Value = 4.397789
$ws.Range("O3").Value = 0.5631392661118858
$ws.Range("P3").Value = 0.5631392661118859
$ws.Range("Q3").Value = 36.89885749112322
$ws.Range("R3").Value = 332.0897174201089
$ws.Range("S3").Value = 0.5631392661118858
$ws.Range("T3").Value = 0.5631392661118859

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.465929666666667
$ws.Range("H4").Value = 4.397789
$ws.Range("M4").Value = 17.866195
$ws.Range("N4").Value = 53.598585
$ws.Range("O4").Value = 0.399712836033041
$ws.Range("P4").Value = 0.399712836033041
$ws.Range("Q4").Value = 26.19058528095167
$ws.Range("R4").Value = 235.715267528565
$ws.Range("S4").Value = 0.399712836033041
$ws.Range("T4").Value = 0.399712836033041
